$d = $word.ActiveDocument

# The three inline pictures in the headers/footers (Pearson logo x2 in the
# footers, BTec logo in the header) have their wp:docPr/pic:cNvPr "name"
# attribute swapped: the Pearson logo images become "image2.png" and the
# BTec logo image becomes "image1.jpg". The "descr" attributes (and the
# actual embedded media/relationship targets) are left untouched.
#
# Word's InlineShape object model does not expose the drawing "name"
# attribute as a settable property (Title maps to the docx "title"
# attribute, AlternativeText maps to "descr"), so perform the edit via a
# WordOpenXML round-trip, which is a lossless way to reach the raw part
# XML through the Word object model.

$xml = $d.WordOpenXML

$xml = $xml.Replace('name="image1.png"', 'name="image2.png"')
$xml = $xml.Replace('name="image2.jpg"', 'name="image1.jpg"')

$d.WordOpenXML = $xml
